# Insert a new data row at row 168 (pushing the existing rows 168-198 down
# to 169-199), then populate the newly-inserted row with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("168:168").Insert()

$ws.Cells.Item(168, 1).Value = 10
$ws.Cells.Item(168, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(168, 3).Value = "La Araucanía"
$ws.Cells.Item(168, 4).Value = 44474
$ws.Cells.Item(168, 5).Value = 9
$ws.Cells.Item(168, 6).Value = 100112009
$ws.Cells.Item(168, 7).Value = "Acelga"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 20
$ws.Cells.Item(168, 11).Value = 9000
$ws.Cells.Item(168, 12).Value = 9000
$ws.Cells.Item(168, 13).Value = 9000
$ws.Cells.Item(168, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(168, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(168, 16).Value = 750
$ws.Cells.Item(168, 17).Value = 12
$ws.Cells.Item(168, 18).Value = "Hortaliza"
